$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching style of existing header row (bold/centered/bordered)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the style from an existing header cell (A1) onto the new header cells,
# without disturbing the values we just set.
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Add new data cells for row 2
$ws.Range("G2").Value = 0.1218615918667638
$ws.Range("H2").Value = 0.991
